$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.469.06"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.840.37"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "260.62"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5252"
$ws.Range("E7").Value = "  +1.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3195"
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06787"
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.77"
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7863"
$ws.Range("E11").Value = "  +3.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07746"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.834.01"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.70"
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.014"
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.85"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007947"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.495.96"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.073.26"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.626"
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.978"
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.375"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.21"
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.176"
$ws.Range("E26").Value = "  -2.32%  "
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.94"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.78"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.157"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08684"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.073"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04869"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7277"
$ws.Range("E34").Value = "  +4.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.857"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.091"
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.247"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01754"
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4777"
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8919"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.51"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.929"
$ws.Range("E43").Value = "  -2.55%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.684"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4173"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.958"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05848"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("E49").Value = "  +1.44%  "
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8910"
$ws.Range("E51").Value = "  +1.50%  "
